# break out stock.yaml completed
# Fix BSE-code column (D) for rows 173-187 on the "day" sheet: these were
# written as text (inlineStr) by mistake, convert them to real numbers.
# Then append the next scraped batch (rows 188-194) of stock data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- Part 1: convert D173:D187 from text to numeric values -----------------
$bseCodes = @{
    173 = 532538
    174 = 539523
    175 = 505200
    176 = 500410
    177 = 532830
    178 = 500483
    179 = 539957
    180 = 500575
    181 = 543066
    182 = 500425
    183 = 531642
    184 = 539336
    185 = 512070
    186 = 533519
    187 = 540750
}

foreach ($r in $bseCodes.Keys) {
    $ws.Cells.Item($r, 4).Value = $bseCodes[$r]
}

# --- Part 2: append new rows 188-194 ---------------------------------------
$newRows = @(
    @(1, "PAGEIND",    "Page Industries Limited",                         "532827", -1.98, 40143.05, 27609,    "day", "18/07/2024 11:35:40"),
    @(2, "MARUTI",     "Maruti Suzuki India Limited",                     "532500",  0.04, 12644.05, 391216,   "day", "18/07/2024 11:35:40"),
    @(3, "DIXON",      "Dixon Technologies",                              "540699", -4.83, 11945.85, 768679,   "day", "18/07/2024 11:35:40"),
    @(4, "BALKRISIND", "Balkrishna Industries Limited",                   "502355",  1.79, 3169.9,   1028346,  "day", "18/07/2024 11:35:40"),
    @(5, "TVSMOTOR",   "Tvs Motor Company Limited",                       "532343", -1.95, 2407.6,   894470,   "day", "18/07/2024 11:35:40"),
    @(6, "HDFCBANK",   "Hdfc Bank Limited",                               "500180", -0.31, 1614.8,   16858532, "day", "18/07/2024 11:35:40"),
    @(7, "M&MFIN",     "Mahindra & Mahindra Financial Services Limited",  "532720", -0.71, 295.2,    2203893,  "day", "18/07/2024 11:35:40")
)

$r = 188
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    # BSE code stored as text (quote-prefixed) to match the source data's
    # inlineStr representation instead of auto-converting to a number.
    $ws.Cells.Item($r, 4).Value = "'" + $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $r = $r + 1
}
